$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Set the value for "Name" property (row 4, column B) which was previously empty
$ws.Range("B4").Value = "TyperoleVs"

# Update the Date value (row 8, column B)
$ws.Range("B8").Value = "2025-07-18T06:40:38+00:00"
